# Apply the "add a 2022 data column" edit to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: columns A-C all become a single uniform width ---
$ws.Range("A1:C1").ColumnWidth = 32.6

# --- New column S: header year + value, cloned (formats + formulas) from column R ---
$ws.Range("R3").Copy($ws.Range("S3"))
$ws.Range("S3").Value = 2022

$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("S4").Value = 0.071025550219041236

# --- Move the active selection ---
$ws.Range("F14").Select()
